# Update sample subscriber data:
#  1. Fix the corrupted city value in row 5 (was "New YorkLos Angeles", should be "Los Angeles")
#  2. Add a new "product_qty" column (J) with a quantity value for every data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix E5 -- the city for Betty L. Profitt should just be "Los Angeles"
$ws.Range("E5").Value = "Los Angeles"

# 2) Add the new product_qty column
$ws.Range("J1").Value = "product_qty"
# Match the left-aligned header style used by the other header cells (A1:I1)
$ws.Range("J1").HorizontalAlignment = -4131

$qty = @{
    2  = 341
    3  = 325
    4  = 425
    5  = 234
    6  = 42
    7  = 23
    8  = 4
    9  = 2345
    10 = 24
    11 = 223
    12 = 76
    13 = 83
}

foreach ($row in $qty.Keys) {
    $ws.Cells.Item($row, 10).Value = $qty[$row]
}

# Move the active selection to reflect where the user ended up after editing
$ws.Range("J14").Select()
